$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row before row 41 (shifts old rows 41-46 down to 42-47),
#    so the period table can grow from 25 data rows (16-40) to 26 (16-41),
#    and the footer rows (signature line / legal-rep labels) end up one
#    row lower than before (46-47 instead of 45-46).
# ---------------------------------------------------------------------------
$ws.Range("B41:J41").Insert(-4121)

# ---------------------------------------------------------------------------
# 2) The newly inserted row 41 should become the new "last row" of the
#    period table (heavier bottom border). Copy that special formatting +
#    values from the old last row (now still row 40) down into row 41,
#    then reset row 40 back to the regular "middle of table" style that
#    rows 16-39 use (copied from row 39).
# ---------------------------------------------------------------------------
$ws.Range("B40:J40").Copy($ws.Range("B41:J41"))
$ws.Range("B39:J39").Copy($ws.Range("B40:J40"))

# ---------------------------------------------------------------------------
# 3) Re-populate the period labels (column E, rows 16-41) in ascending
#    chronological order (was descending, newest-first before) and add the
#    new period 2508 in the new row 41.
# ---------------------------------------------------------------------------
$periods = @("2307","2308","2309","2310","2311","2312", `
             "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412", `
             "2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# ---------------------------------------------------------------------------
# 4) Update the summary figures at the top of the statement:
#    - VALOR MORA total (E11): 1160000 -> 1206400
#    - Cant. Periodos (F13): 25 -> 26 (one more period now tracked)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1206400
$ws.Range("F13").Value = 26

Write-Output "edit applied"
